$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.07%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-0.51%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.254"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.47%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'0.33%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.40%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9181"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.74%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.444"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.81%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'14.15%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1834"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.25%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09245"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.45%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04263"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.82%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.20%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.95%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005745"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.96%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E17").Value = "'-0.12%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.315"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.26%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'7.416"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'13.27%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1384"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.74%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2894"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.87%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04074"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.78%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.32%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004295"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.44%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001273"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.11%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02467"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.11%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05282"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.86%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007855"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.25%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'1.12%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006820"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.00%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-1.85%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.14%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = "'0.00006729"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.91%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.20%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2055"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,849.99%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-2.43%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.20%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.20%"
$ws.Range("E51").Style = "Normal"
